$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing forecast-error values for rows 2-10 (Q0..Q8) ---
$ws.Range("B2").Value = 0.001068801475677701
$ws.Range("C2").Value = 0.5121987664681366
$ws.Range("D2").Value = 0.5908018752482725
$ws.Range("E2").Value = 0.7686363738779688
$ws.Range("F2").Value = 0.776283934891037
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = 0.1076112683830704
$ws.Range("C3").Value = 0.5501073581205883
$ws.Range("D3").Value = 0.7154227106871884
$ws.Range("E3").Value = 0.8458266434011099
$ws.Range("F3").Value = 0.8474707421379045
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = 0.01447167482266423
$ws.Range("C4").Value = 0.6025716605263358
$ws.Range("D4").Value = 0.7372352437580105
$ws.Range("E4").Value = 0.8586240409853492
$ws.Range("F4").Value = 0.8673987081222851
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = 0.107527968228121
$ws.Range("C5").Value = 0.6082289609806715
$ws.Range("D5").Value = 0.7759918099762468
$ws.Range("E5").Value = 0.8809039731867753
$ws.Range("F5").Value = 0.883568901926184
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = 0.01941581067445748
$ws.Range("C6").Value = 0.5766303201776867
$ws.Range("D6").Value = 0.6540589049982858
$ws.Range("E6").Value = 0.8087390833873962
$ws.Range("F6").Value = 0.8172468464677016
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = 0.1112754059595625
$ws.Range("C7").Value = 0.564041488129664
$ws.Range("D7").Value = 0.7222269590935267
$ws.Range("E7").Value = 0.8498393725249065
$ws.Range("F7").Value = 0.8518327715854772
$ws.Range("G7").Value = 46

$ws.Range("B8").Value = 0.01597012535789634
$ws.Range("C8").Value = 0.5344976713245601
$ws.Range("D8").Value = 0.6191284896479551
$ws.Range("E8").Value = 0.7868471831607171
$ws.Range("F8").Value = 0.7955744781642835
$ws.Range("G8").Value = 45

$ws.Range("B9").Value = 0.05572434282723118
$ws.Range("C9").Value = 0.6058820781989243
$ws.Range("D9").Value = 0.7089638270491668
$ws.Range("E9").Value = 0.8419998972975987
$ws.Range("F9").Value = 0.849867014019446
$ws.Range("G9").Value = 44

$ws.Range("B10").Value = 0.06564388546255764
$ws.Range("C10").Value = 0.6068818695593868
$ws.Range("D10").Value = 0.7068678503580751
$ws.Range("E10").Value = 0.8407543341298188
$ws.Range("F10").Value = 0.8481074928832305
$ws.Range("G10").Value = 43

# --- Add new row 11 (Q9) reusing the same formatting as row 10's label cell ---
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.03898798267175772
$ws.Range("C11").Value = 0.5674245104323993
$ws.Range("D11").Value = 0.5443824653879992
$ws.Range("E11").Value = 0.7378227872517893
$ws.Range("F11").Value = 0.7457231047781384
$ws.Range("G11").Value = 42
